# Update the Tasks sheet: rename existing tasks, give Task 2 a status,
# and append placeholder Task 3 .. Task 10 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the order in which new string values are first assigned controls
# the order they are appended to the shared-string table, so we write the
# "Task 2" row before the "Task 1" row to match the target ordering.
$ws.Range("A3").Value = "Task 2: Home Page that auto adjusts on zoom levels"
$ws.Range("B3").Value = "In Development"
$ws.Range("A2").Value = "Task 1: Layout including footer and header that auto adjusts on zoom levels"

$ws.Range("A4").Value = "Task 3:"
$ws.Range("A5").Value = "Task 4:"
$ws.Range("A6").Value = "Task 5:"
$ws.Range("A7").Value = "Task 6:"
$ws.Range("A8").Value = "Task 7:"
$ws.Range("A9").Value = "Task 8:"
$ws.Range("A10").Value = "Task 9:"
$ws.Range("A11").Value = "Task 10:"

# Match the saved selection in the target workbook.
$ws.Range("A4").Select()
